$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("D2").Value = 0.0001334558586972819
$ws.Range("E2").Value = 0.0001334558586972819
$ws.Range("D3").Value = 0.9714458571676261
$ws.Range("E3").Value = 0.9714458571676261
$ws.Range("D4").Value = 0.000023425184784327210158350746072031256517220754176378250122
$ws.Range("E4").Value = 0.000023425184784327210158350746072031256517220754176378250122
$ws.Range("D5").Value = 0.000000000000000047925567663770038300149908562304766242674177
$ws.Range("E5").Value = 0.000000000000000047925567663770038300149908562304766242674177
$ws.Range("D6").Value = 0.7087541795104007
$ws.Range("E6").Value = 0.7087541795104007
$ws.Range("D7").Value = 0.9976365537575965
$ws.Range("E7").Value = 0.002363446242403477
$ws.Range("D8").Value = 0.8960762728884814
$ws.Range("E8").Value = 0.1039237271115186
$ws.Range("D9").Value = 0.9668134181145299
$ws.Range("E9").Value = 0.03318658188547008
$ws.Range("D10").Value = 0.997682432714646
$ws.Range("E10").Value = 0.002317567285354039
$ws.Range("D11").Value = 0.9793656142634262
$ws.Range("E11").Value = 0.02063438573657383
$ws.Range("F11").Value = 0.4958714842796326
$ws.Range("D12").Value = 0.0002022080920420974
$ws.Range("E12").Value = 0.0002022080920420974
$ws.Range("D13").Value = 0.9964219922605527
$ws.Range("E13").Value = 0.9964219922605527
$ws.Range("D14").Value = 0.0004056613643315822
$ws.Range("E14").Value = 0.0004056613643315822
$ws.Range("D15").Value = 0.000000000000000000007636477246211202312497242724655315698745
$ws.Range("E15").Value = 0.000000000000000000007636477246211202312497242724655315698745
$ws.Range("D16").Value = 0.1322650787433982
$ws.Range("E16").Value = 0.1322650787433982
$ws.Range("D17").Value = 0.9989313939047907
$ws.Range("E17").Value = 0.001068606095209312
$ws.Range("D18").Value = 0.9948670563109221
$ws.Range("E18").Value = 0.005132943689077907
$ws.Range("D19").Value = 0.991644591638772
$ws.Range("E19").Value = 0.008355408361227945
$ws.Range("D20").Value = 0.9983848385020777
$ws.Range("E20").Value = 0.001615161497922291
$ws.Range("D21").Value = 0.9831201494415138
$ws.Range("E21").Value = 0.01687985055848618
$ws.Range("F21").Value = 0.5808672308921814
